$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# hunk @ old line 727
$ws.Range("H2").Value = 465
$ws.Range("I2").Value = 465
$ws.Range("K2").Value = 465
$ws.Range("M2").Value = -352

# hunk @ old line 877
$ws.Range("H5").Value = 124
$ws.Range("I5").Value = 124
$ws.Range("K5").Value = 124
$ws.Range("M5").Value = -9

# hunk @ old line 2206
$ws.Range("I32").Value = 3995.6667
$ws.Range("J32").Value = 3500
$ws.Range("K32").Value = 3995.6667
$ws.Range("L32").Value = 3500
$ws.Range("M32").Value = -3669.6667
$ws.Range("N32").Value = -4152

# hunk @ old line 3137
$ws.Range("H51").Value = 7999.6665
$ws.Range("J51").Value = 7999.6665
$ws.Range("L51").Value = 7999.6665
$ws.Range("N51").Value = -8967.666499999999

# hunk @ old line 3774
$ws.Range("H64").Value = 3467
$ws.Range("J64").Value = 3520.6
$ws.Range("L64").Value = 3520.6
$ws.Range("N64").Value = -4016.6

# hunk @ old line 3924
$ws.Range("H67").Value = 3467
$ws.Range("J67").Value = 3520.6
$ws.Range("L67").Value = 3520.6
$ws.Range("N67").Value = -5236.6

# hunk @ old line 4270
$ws.Range("H74").Value = 5428.5713
$ws.Range("I74").Value = 5428.5713
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 5428.5713
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4492.5713
$ws.Range("N74").ClearContents()

# hunk @ old line 4368
$ws.Range("H76").Value = 5704.7393
$ws.Range("I76").Value = 4247.273
$ws.Range("K76").Value = 4247.273
$ws.Range("M76").Value = -3932.273

# hunk @ old line 4420
$ws.Range("H77").Value = 5428.5713
$ws.Range("I77").Value = 5428.5713
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 27142.8565
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -22462.8565
$ws.Range("N77").ClearContents()

# hunk @ old line 4518
$ws.Range("H79").Value = 5704.7393
$ws.Range("I79").Value = 4247.273
$ws.Range("K79").Value = 4247.273
$ws.Range("M79").Value = -3155.273

# hunk @ old line 5164
$ws.Range("H92").Value = 457.89474
$ws.Range("I92").Value = 332.73334
$ws.Range("J92").Value = 927.25
$ws.Range("K92").Value = 332.73334
$ws.Range("L92").Value = 927.25
$ws.Range("M92").Value = 915.26666
$ws.Range("N92").Value = -3423.25

# hunk @ old line 5363
$ws.Range("H96").Value = 2387.4546
$ws.Range("I96").Value = 1177.5714
$ws.Range("K96").Value = 3532.7142
$ws.Range("M96").Value = -2159.7142

# hunk @ old line 5568
$ws.Range("H100").Value = 1319.8667
$ws.Range("J100").Value = 2046.6
$ws.Range("L100").Value = 2046.6
$ws.Range("N100").Value = -3128.6

# hunk @ old line 6790
$ws.Range("H125").Value = 989.25
$ws.Range("I125").Value = 986
$ws.Range("K125").Value = 8874
$ws.Range("M125").Value = -6414

# hunk @ old line 7081
$ws.Range("H131").Value = 4396.4
$ws.Range("I131").Value = 1577.3334
$ws.Range("J131").Value = 8625
$ws.Range("K131").Value = 4732.0002
$ws.Range("L131").Value = 25875
$ws.Range("M131").Value = 307.9997999999996
$ws.Range("N131").Value = -35955

# hunk @ old line 7133
$ws.Range("H132").Value = 1239.7587
$ws.Range("I132").Value = 1273.3214
$ws.Range("K132").Value = 3819.9642
$ws.Range("M132").Value = -1289.9642

$ws = $wb.Worksheets.Item("ARM")
# hunk @ old line 10596
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# hunk @ old line 14196
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# hunk @ old line 19014
$ws.Range("H94").Value = 1837.75
$ws.Range("I94").Value = 2042.5
$ws.Range("J94").Value = 404.5
$ws.Range("K94").Value = 2042.5
$ws.Range("L94").Value = 404.5
$ws.Range("M94").Value = -1591.5
$ws.Range("N94").Value = -1306.5

# hunk @ old line 19259
$ws.Range("H99").Value = 3290.923
$ws.Range("J99").Value = 3295.25
$ws.Range("L99").Value = 3295.25
$ws.Range("N99").Value = -6291.25

# hunk @ old line 20923
$ws.Range("H134").Value = 1955.5
$ws.Range("I134").Value = 1279.0834
$ws.Range("K134").Value = 3837.2502
$ws.Range("M134").Value = -1302.2502

$ws = $wb.Worksheets.Item("CRP")
# hunk @ old line 22383
$ws.Range("H22").Value = 589.5454999999999
$ws.Range("I22").Value = 407.8
$ws.Range("J22").Value = 741
$ws.Range("K22").Value = 407.8
$ws.Range("L22").Value = 741
$ws.Range("M22").Value = -57.80000000000001
$ws.Range("N22").Value = -1441

# hunk @ old line 27253
$ws.Range("H122").Value = 4831.6
$ws.Range("I122").Value = 4641
$ws.Range("K122").Value = 13923
$ws.Range("M122").Value = -11473

# hunk @ old line 27734
$ws.Range("H132").Value = 3922
$ws.Range("J132").Value = 5894.25
$ws.Range("L132").Value = 17682.75
$ws.Range("N132").Value = -22742.75

$ws = $wb.Worksheets.Item("CUL")
# hunk @ old line 30271
$ws.Range("H41").Value = 347.5
$ws.Range("I41").Value = 233.33333
$ws.Range("J41").Value = 690
$ws.Range("K41").Value = 699.99999
$ws.Range("L41").Value = 2070
$ws.Range("M41").Value = -361.99999
$ws.Range("N41").Value = -2746

# hunk @ old line 32479
$ws.Range("H86").Value = 92.71429000000001
$ws.Range("I86").Value = 91.666664
$ws.Range("J86").Value = 99
$ws.Range("K86").Value = 274.999992
$ws.Range("L86").Value = 297
$ws.Range("M86").Value = 911.000008
$ws.Range("N86").Value = -2669

# hunk @ old line 32626
$ws.Range("H89").Value = 92.71429000000001
$ws.Range("I89").Value = 91.666664
$ws.Range("J89").Value = 99
$ws.Range("K89").Value = 824.9999759999999
$ws.Range("L89").Value = 891
$ws.Range("M89").Value = 5103.000024
$ws.Range("N89").Value = -12747

$ws = $wb.Worksheets.Item("GSM")
# hunk @ old line 37004
$ws.Range("H36").Value = 11000
$ws.Range("I36").Value = 14000
$ws.Range("J36").Value = 8000
$ws.Range("K36").Value = 14000
$ws.Range("L36").Value = 8000
$ws.Range("M36").Value = -13515
$ws.Range("N36").Value = -8970

# hunk @ old line 39820
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492

# hunk @ old line 40166
$ws.Range("H102").Value = 1264.375
$ws.Range("I102").Value = 296.7647
$ws.Range("K102").Value = 296.7647
$ws.Range("M102").Value = 1325.2353

# hunk @ old line 41609
$ws.Range("H132").Value = 1921.9412
$ws.Range("I132").Value = 760.3570999999999
$ws.Range("J132").Value = 7342.6665
$ws.Range("K132").Value = 2281.0713
$ws.Range("L132").Value = 22027.9995
$ws.Range("M132").Value = 248.9287000000004
$ws.Range("N132").Value = -27087.9995

$ws = $wb.Worksheets.Item("LTW")
# hunk @ old line 46577
$ws.Range("H93").Value = 997.5
$ws.Range("I93").Value = 997.5
$ws.Range("K93").Value = 997.5
$ws.Range("M93").Value = 250.5

# hunk @ old line 48440
$ws.Range("H132").Value = 3604.2307
$ws.Range("I132").Value = 2510.1875
$ws.Range("J132").Value = 5354.7
$ws.Range("K132").Value = 7530.5625
$ws.Range("L132").Value = 16064.1
$ws.Range("M132").Value = -5000.5625
$ws.Range("N132").Value = -21124.1

# hunk @ old line 48633
$ws.Range("H136").Value = 2397.5
$ws.Range("I136").Value = 2397.5
$ws.Range("K136").Value = 7192.5
$ws.Range("M136").Value = -4642.5
